# Commit: "simplify nixie daughter board"
# Remove cathode pins that would not be usable for current sensing anyway.
# Pivot to using anode pin as shunt measurement with a R2R flash ADC (not yet implemented).
#
# This adds new "resistor shunt" calcs (voltage drop & power dissipation across
# the existing shunt resistor R, reusing Imin/Inom/Imax already on the sheet)
# to the "Nixie limiting" worksheet.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# New labels in column A. Order matters: it controls the order new strings are
# appended to the shared-string table, so Vres-max is written first.
$ws2.Range("A10").Value = "Vres-max"
$ws2.Range("A8").Value  = "Vres-min"
$ws2.Range("A9").Value  = "Vres-nom"
$ws2.Range("A12").Value = "Pres-min"
$ws2.Range("A13").Value = "Pres-nom"
$ws2.Range("A14").Value = "Pres-max"

# Units column. "V" already exists in the shared-string table; "W" is new and
# must be introduced here (before it is reused on rows 13/14).
$ws2.Range("C12").Value = "W"
$ws2.Range("C8").Value  = "V"
$ws2.Range("C9").Value  = "V"
$ws2.Range("C10").Value = "V"
$ws2.Range("C13").Value = "W"
$ws2.Range("C14").Value = "W"

# Voltage across the shunt resistor for min/nom/max current, using the
# existing Imin/Inom/Imax (B4/B5/B6) and shunt resistance R (B3).
$ws2.Range("B8").Formula  = "=B4*B3"
$ws2.Range("B9").Formula  = "=B5*B3"
$ws2.Range("B10").Formula = "=B6*B3"

# Power dissipated in the shunt resistor for min/nom/max current.
$ws2.Range("B12").Formula = "=B4*B8"
$ws2.Range("B13").Formula = "=B5*B9"
$ws2.Range("B14").Formula = "=B6*B10"

# Match the scientific-notation number format already used by B4:B6.
$fmt = $ws2.Range("B4").NumberFormat
$ws2.Range("B8").NumberFormat  = $fmt
$ws2.Range("B9").NumberFormat  = $fmt
$ws2.Range("B10").NumberFormat = $fmt
$ws2.Range("B12").NumberFormat = $fmt
$ws2.Range("B13").NumberFormat = $fmt
$ws2.Range("B14").NumberFormat = $fmt

# Column A widens slightly to fit the new "Vres-xxx"/"Pres-xxx" labels.
$ws2.Columns.Item(1).ColumnWidth = 8.67
